$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row for 2025-01-23 (serial 45680) - becomes the new last row
$ws.Cells.Item(93, 1).Value = 45680
$ws.Cells.Item(93, 2).Value = 220
$ws.Cells.Item(93, 3).Value = 218
$ws.Cells.Item(93, 4).Value = 216

# The "last row" date cell (A92) carries a special highlight style; move that
# style onto the new last row (A93) first...
$ws.Range("A92").Copy() | Out-Null
$ws.Range("A93").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# ...then revert A92 back to the normal date style shared by the other rows.
$ws.Range("A91").Copy() | Out-Null
$ws.Range("A92").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0
